$d = $word.ActiveDocument

# The document currently ends with a paragraph containing "+3 more".
# We replace that with the new team member block:
#   Calistas Mlilo
#   cmlilo2004@gmail.com        (hyperlink)
#   +27 67 820 6902
#   <empty paragraph>
#   Jabulile Msibi:
#   jabulilemsibi184@gmail.com  (hyperlink)
#   +27 67 838 4308

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Text = "Calistas Mlilo`rcmlilo2004@gmail.com`r+27 67 820 6902`r`rJabulile Msibi: `rjabulilemsibi184@gmail.com`r+27 67 838 4308"

# Re-fetch paragraph count/indices after the multi-paragraph insert.
$total = $d.Paragraphs.Count

# Turn the "cmlilo2004@gmail.com" paragraph into a mailto hyperlink.
$emailPara1 = $d.Paragraphs.Item($total - 5)
$emailRange1 = $emailPara1.Range
[void]$emailRange1.MoveEnd(1, -1)
$d.Hyperlinks.Add($emailRange1, "mailto:cmlilo2004@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "cmlilo2004@gmail.com") | Out-Null

# Turn the "jabulilemsibi184@gmail.com" paragraph into a mailto hyperlink.
$emailPara2 = $d.Paragraphs.Item($total - 1)
$emailRange2 = $emailPara2.Range
[void]$emailRange2.MoveEnd(1, -1)
$d.Hyperlinks.Add($emailRange2, "mailto:jabulilemsibi184@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "jabulilemsibi184@gmail.com") | Out-Null
